# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same list of events (mirrored data), but the
# scraped counters drifted slightly differently between the two exports,
# so each sheet gets its own explicit set of F-column values.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# row -> new value for "展览"
$exhibitUpdates = @{
    2  = 1062
    3  = 754
    4  = 260
    8  = 1687
    9  = 6394
    10 = 484
    11 = 367
    12 = 308
    13 = 0
    14 = 378
    15 = 138
    16 = 6375
    17 = 275
    18 = 1288
    19 = 138
    20 = 117
    21 = 221
    22 = 108
    23 = 275
    24 = 0
    26 = 10
    27 = 98
    28 = 10
    29 = 391
    30 = 92
    33 = 48
    35 = 25
    36 = 63
    37 = 63
}

# row -> new value for "全部类型"
$allUpdates = @{
    2  = 1062
    3  = 754
    4  = 260
    8  = 1687
    9  = 6394
    10 = 484
    11 = 367
    12 = 308
    13 = 102
    14 = 378
    15 = 138
    16 = 6375
    17 = 275
    18 = 1288
    19 = 138
    20 = 117
    21 = 221
    22 = 108
    23 = 275
    24 = 106
    26 = 10
    27 = 98
    28 = 10
    29 = 391
    30 = 92
    33 = 48
    35 = 25
    36 = 63
    37 = 0
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
